$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Flip three existing test rows from "YES" to "NO" in the Execution Flag column (C)
$ws.Range("C2").Value = "NO"
$ws.Range("C3").Value = "NO"
$ws.Range("C13").Value = "NO"

# Append a new test-case row (row 19) describing a checkout-page check
$ws.Range("A19").Value = "checking checkout page"
$ws.Range("B19").Value = ([string][char]0x201C) + ([string][char]0x201D)
$ws.Range("C19").Value = "YES"
$ws.Range("D19").Value = "login"
$ws.Range("E19").Value = "checkout"

# Move the active selection to C11, matching the author's saved cursor position
$ws.Range("C11").Select()
